# Re-pull / push updated "dSF" (column F) data for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = -1
    13 = -2
    19 = -1
    24 = -3
    26 = 7
    32 = -1
    35 = 0
    49 = -3
    56 = 0
    58 = -1
    59 = -1
    61 = -1
    62 = -1
    65 = 2
    67 = -1
    68 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
